$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2032520325203252
$ws.Range("C2").Value = 0.5772357723577236
$ws.Range("P2").Value = 0.1585365853658537
$ws.Range("S2").Value = 0.06097560975609756
$ws.Range("B3").Value = 0.006493506493506494
$ws.Range("C3").Value = 0.05844155844155844
$ws.Range("J3").Value = 0.01298701298701299
$ws.Range("P3").Value = 0.7337662337662337
$ws.Range("S3").Value = 0.1883116883116883
$ws.Range("J4").Value = 0.02564102564102564
$ws.Range("P4").Value = 0.5897435897435898
$ws.Range("S4").Value = 0.3846153846153846
$ws.Range("B6").Value = 0.07936507936507936
$ws.Range("F6").Value = 0.06878306878306878
$ws.Range("J6").Value = 0.2222222222222222
$ws.Range("O6").Value = 0.01587301587301587
$ws.Range("Q6").Value = 0.1375661375661376
$ws.Range("S6").Value = 0.3650793650793651
$ws.Range("B7").Value = 0.1104294478527607
$ws.Range("D7").Value = 0.01840490797546012
$ws.Range("E7").Value = 0.006134969325153374
$ws.Range("F7").Value = 0.049079754601227
$ws.Range("J7").Value = 0.09202453987730061
$ws.Range("O7").Value = 0.01840490797546012
$ws.Range("Q7").Value = 0.1533742331288344
$ws.Range("R7").Value = 0.0736196319018405
$ws.Range("S7").Value = 0.4785276073619632
$ws.Range("B8").Value = 0.08585858585858586
$ws.Range("D8").Value = 0.01262626262626263
$ws.Range("F8").Value = 0.05303030303030303
$ws.Range("J8").Value = 0.1035353535353535
$ws.Range("O8").Value = 0.03535353535353535
$ws.Range("Q8").Value = 0.196969696969697
$ws.Range("R8").Value = 0.08585858585858586
$ws.Range("S8").Value = 0.4267676767676767
$ws.Range("B9").Value = 0.05116279069767442
$ws.Range("D9").Value = 0.02325581395348837
$ws.Range("F9").Value = 0.05116279069767442
$ws.Range("J9").Value = 0.1302325581395349
$ws.Range("O9").Value = 0.02325581395348837
$ws.Range("Q9").Value = 0.1627906976744186
$ws.Range("R9").Value = 0.1023255813953488
$ws.Range("S9").Value = 0.4558139534883721
$ws.Range("B10").Value = 0.09656301145662848
$ws.Range("D10").Value = 0.0220949263502455
$ws.Range("F10").Value = 0.0630114566284779
$ws.Range("J10").Value = 0.1170212765957447
$ws.Range("O10").Value = 0.01554828150572831
$ws.Range("Q10").Value = 0.2356792144026187
$ws.Range("R10").Value = 0.102291325695581
$ws.Range("S10").Value = 0.3477905073649755
$ws.Range("G11").Value = 0.1176470588235294
$ws.Range("J11").Value = 0.08403361344537816
$ws.Range("K11").Value = 0.1764705882352941
$ws.Range("L11").Value = 0.5966386554621849
$ws.Range("S11").Value = 0.02521008403361345
$ws.Range("G12").Value = 0.7046979865771812
$ws.Range("J12").Value = 0.2214765100671141
$ws.Range("K12").Value = 0.01342281879194631
$ws.Range("L12").Value = 0.03355704697986577
$ws.Range("S12").Value = 0.02684563758389262
$ws.Range("F13").Value = 0.02
$ws.Range("G13").Value = 0.7
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.08
$ws.Range("F15").Value = 0.01694915254237288
$ws.Range("H15").Value = 0.1398305084745763
$ws.Range("I15").Value = 0.06779661016949153
$ws.Range("J15").Value = 0.4364406779661017
$ws.Range("K15").Value = 0.04661016949152542
$ws.Range("M15").Value = 0.008474576271186441
$ws.Range("O15").Value = 0.07627118644067797
$ws.Range("S15").Value = 0.2076271186440678
$ws.Range("F16").Value = 0.01169590643274854
$ws.Range("H16").Value = 0.1812865497076023
$ws.Range("I16").Value = 0.06432748538011696
$ws.Range("J16").Value = 0.4385964912280702
$ws.Range("K16").Value = 0.06432748538011696
$ws.Range("M16").Value = 0.03508771929824561
$ws.Range("O16").Value = 0.04678362573099415
$ws.Range("S16").Value = 0.1578947368421053
$ws.Range("F17").Value = 0.02022471910112359
$ws.Range("H17").Value = 0.1662921348314607
$ws.Range("I17").Value = 0.1213483146067416
$ws.Range("J17").Value = 0.4067415730337079
$ws.Range("K17").Value = 0.09662921348314607
$ws.Range("M17").Value = 0.01573033707865169
$ws.Range("O17").Value = 0.07191011235955057
$ws.Range("S17").Value = 0.101123595505618
$ws.Range("F18").Value = 0.004672897196261682
$ws.Range("H18").Value = 0.1775700934579439
$ws.Range("I18").Value = 0.1121495327102804
$ws.Range("J18").Value = 0.4672897196261682
$ws.Range("K18").Value = 0.06542056074766354
$ws.Range("M18").Value = 0.02803738317757009
$ws.Range("O18").Value = 0.06074766355140187
$ws.Range("S18").Value = 0.08411214953271028
$ws.Range("F19").Value = 0.01456726649528706
$ws.Range("H19").Value = 0.194515852613539
$ws.Range("I19").Value = 0.09511568123393316
$ws.Range("J19").Value = 0.3787489288774636
$ws.Range("K19").Value = 0.09682947729220223
$ws.Range("M19").Value = 0.02656383890317052
$ws.Range("N19").Value = 0.000856898029134533
$ws.Range("O19").Value = 0.0805484147386461
$ws.Range("S19").Value = 0.1122536418166238
